$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, shifting existing rows 104:110 down to 105:111.
$ws.Rows.Item(104).Insert()

# Populate the new row 104 with the new weekly price observation.
$ws.Cells.Item(104, 1).Value = 10
$ws.Cells.Item(104, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(104, 3).Value = "La Araucanía"
$ws.Cells.Item(104, 4).Value = 45132
$ws.Cells.Item(104, 5).Value = 9
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100108
$ws.Cells.Item(104, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(104, 9).Value = 100108004
$ws.Cells.Item(104, 10).Value = "Papaya"
$ws.Cells.Item(104, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(104, 12).Value = "Primera"
$ws.Cells.Item(104, 13).Value = 40
$ws.Cells.Item(104, 14).Value = 25000
$ws.Cells.Item(104, 15).Value = 25000
$ws.Cells.Item(104, 16).Value = 25000
$ws.Cells.Item(104, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(104, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(104, 19).Value = 2500
$ws.Cells.Item(104, 20).Value = 10
